$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.335339188575745
$ws.Range("B1").Value = 1.59981906414032
$ws.Range("C1").Value = 4.048702239990234
$ws.Range("D1").Value = 3.23529314994812
$ws.Range("E1").Value = 1.102744460105896
